# Apply the "Enemies shoot now, but it's in progress" edit to the
# "Towers vs. Enemies chart" worksheet (Sheet1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Header / label text tweaks -------------------------------------------------
$ws.Range("C23").Value = "Damage"            # was "Damage (per sec)"
$ws.Range("H29").Value = "very slow"         # was "slow"
$ws.Range("B34").Value = "Reload"            # was "Speed"
$ws.Range("C34").Value = "Damage"            # was "Damage(per sec)"

# --- Towers table (rows 24-32) numeric tweaks -----------------------------------
$ws.Range("F24").Value = 2                   # was 0
$ws.Range("C27").Value = 8                   # was 12
$ws.Range("C28").Value = 6                   # was 15
$ws.Range("C29").Value = 1000                # was 35
$ws.Range("C31").Value = 25                  # was 5

# --- Enemies table (rows 35-45) numeric tweaks ----------------------------------
$ws.Range("B35").Value = 25                  # was 2
$ws.Range("C35").Value = 2                   # was 3
$ws.Range("B36").Value = 100                 # was 1
$ws.Range("C36").Value = 10                  # was 5
$ws.Range("B37").Value = 500                 # was 3
$ws.Range("C37").Value = 50                  # was 2
$ws.Range("B38").Value = 150                 # was 3
$ws.Range("C38").Value = 15                  # was 3
$ws.Range("B39").Value = 50                  # was 1
$ws.Range("B40").Value = 25                  # was 2
$ws.Range("B41").Value = 50                  # was "3 or 4" (text)
$ws.Range("C41").Value = 5                   # was 2
$ws.Range("B42").Value = 25                  # was 2
$ws.Range("C42").Value = 5                   # was 3
$ws.Range("B43").Value = 5                   # was 1
$ws.Range("B44").Value = 50                  # was 2
$ws.Range("C44").Value = 5                   # was 3
$ws.Range("B45").Value = 25                  # was 2

# --- View state: scroll position & selected cell --------------------------------
$ws.Activate()
$appWin = $excel.ActiveWindow
$appWin.ScrollRow = 13
$appWin.ScrollColumn = 1
$ws.Range("B45").Select()
